$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log row appended at row 6 (the sheet previously had data through row 5).
# Values must stay as literal text (the source column B holds a plain
# "YYYY-MM-DD" string, not a date), so force Text format before assigning,
# then drop back to the default "Normal" style so no style index is stamped
# on the new cells (matching the unstyled cells used by the other data rows).
$row = 6

$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025-11-01 03:55:45"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").NumberFormat = "@"
$ws.Range("B$row").Value = "2025-10-30"
$ws.Range("B$row").Style = "Normal"

$ws.Range("C$row").NumberFormat = "@"
$ws.Range("C$row").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice30102025-7.pdf"
$ws.Range("C$row").Style = "Normal"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "/home/runner/work/rashtriyametal_downloader/rashtriyametal_downloader/data/RashtriyaMetal/PDFs/ListPrice30102025-7.pdf"
$ws.Range("D$row").Style = "Normal"
